$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 observation station was re-classified from a USACE gauge (01440)
# to an MS River / USGS gauge (291929089562600).
$ws.Range("A12").Value = "USGS"
$ws.Range("B12").Value = "291929089562600"

# Column B now needs to be a bit wider to comfortably show the longer USGS
# station id, so give it a fitted custom width.
$ws.Columns.Item(2).ColumnWidth = 9.14

# Reflect the last selected cell in the saved view.
$ws.Range("F12").Select()
